$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "Jabeessaa Caalaa" (originally row 4)
$ws.Rows.Item(4).Delete()

# After the above deletion, "Baqqalaa Barsiisaa" (originally row 9) is now row 8
$ws.Rows.Item(8).Delete()

# Update the selection to match the target state
$ws.Range("B15").Select()
